$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: split "À chaque visite, le membre doit se présenter..." so that
# " ou professionnel" becomes its own run, inserted right after "le membre".
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("le membre", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'le membre' anchor text"
}
$insertAt = $d.Range($rng1.End, $rng1.End)
$insertAt.InsertAfter(" ou professionnel")

# Force the newly inserted text to live in its own run (distinct from the
# run that follows it) by briefly bookmarking it and removing the bookmark;
# this mirrors how Word naturally splits runs on formatting boundaries
# without leaving stray empty <w:rPr/> behind.
$splitRng1 = $d.Range($rng1.End, $rng1.End + 17)
$d.Bookmarks.Add("TempSplit1", $splitRng1)
$d.Bookmarks("TempSplit1").Delete()

# ---------------------------------------------------------------------------
# Edit 2: split "apparait" into "appar" / "ait" and drop a _GoBack bookmark
# (zero-length) right at that split point. Word keeps only a single
# "_GoBack" bookmark at a time, so adding it here automatically relocates
# (removes) the one that currently sits in the "Accéder aux services..."
# paragraph.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("apparait", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'apparait'"
}
$splitPoint2 = $rng2.Start + 5
$zeroRng2 = $d.Range($splitPoint2, $splitPoint2)
$d.Bookmarks.Add("_GoBack", $zeroRng2)

# ---------------------------------------------------------------------------
# Edit 3: the "Accéder aux services libres..." paragraph now contains the
# same sentence split across two runs (since the _GoBack bookmark moved
# away from there) -- merge them back into a single run, matching the
# target OOXML.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Accéder aux services libres (aucune inscription requise) : Il s’agit uniquement d’authentifier l’identité du membre.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) {
    throw "Could not find the 'Accéder aux services libres' sentence"
}
$mergedText = $rng3.Text
$mergeStart = $rng3.Start
$rng3.Delete()
$mergeInsertAt = $d.Range($mergeStart, $mergeStart)
$mergeInsertAt.InsertAfter($mergedText)

Write-Output "Done"
